$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.325.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.720.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4715"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2627"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06207"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.717.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07071"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5908"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.397"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.322.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006782"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.939.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.550"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.326"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "108.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.406"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.002"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.686"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07735"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04460"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.616"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9742"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6186"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.416"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.11%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01476"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.354"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3812"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1165"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.264"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05288"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.702"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3375"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9230"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("B38").Value = "Quant"
$ws.Range("C38").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "114.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.79%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.25%  "
